# Update the worksheet's division problems to the new generated set.
#
# Row 1 of the table gets two new problems prepended (76÷8=, 10÷8=) and,
# because the table is a fixed 5-column grid, the trailing two problems
# (47÷7=, 61÷6=) fall off the end while the rest shift right by two slots.
# Net effect on row 1, cell-by-cell:
#   51÷6= -> 76÷8=
#   40÷6= -> 10÷8=
#   10÷4= -> 51÷6=
#   47÷7= -> 43÷2=
#   61÷6= -> 34÷4=
# Addressing cells directly (instead of a text Find/Replace) avoids the
# collision created by "51÷6=" being both an old and a new value in this row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "76÷8="
$t.Cell(1, 2).Range.Text = "10÷8="
$t.Cell(1, 3).Range.Text = "51÷6="
$t.Cell(1, 4).Range.Text = "43÷2="
$t.Cell(1, 5).Range.Text = "34÷4="

# The remaining problems (rows 5, 9, 13, 17) are simple 1-for-1 text swaps;
# every old value below is unique across the document.
$replacements = @(
    @("73÷6=", "59÷3="),
    @("30÷4=", "43÷7="),
    @("57÷9=", "57÷4="),
    @("29÷6=", "45÷2="),
    @("63÷8=", "44÷9="),
    @("88÷6=", "89÷6="),
    @("53÷4=", "20÷9="),
    @("37÷8=", "49÷3="),
    @("43÷4=", "64÷4="),
    @("25÷4=", "56÷9="),
    @("37÷6=", "18÷6="),
    @("72÷2=", "90÷8="),
    @("76÷9=", "44÷6="),
    @("22÷3=", "53÷2="),
    @("17÷6=", "21÷3="),
    @("11÷2=", "44÷8="),
    @("96÷2=", "83÷9="),
    @("85÷6=", "71÷6="),
    @("75÷9=", "62÷2="),
    @("29÷8=", "67÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Host "edit applied"
